$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -22.04450000000001
$ws.Range("A21").Value = -20.07989999999998
$ws.Range("A23").Value = -20.23579999999998
$ws.Range("A25").Value = -21.71489999999999
$ws.Range("B27").Value = 5.990200000000001
$ws.Range("B31").Value = 5.433900000000001
$ws.Range("B39").Value = 9.668200000000006
$ws.Range("B48").Value = 5.322200000000004
$ws.Range("B51").Value = 5.549599999999998
$ws.Range("B52").Value = 5.167499999999999
$ws.Range("A53").Value = -21.86360000000001
$ws.Range("B55").Value = 5.952899999999995
$ws.Range("B56").Value = 4.921099999999999
$ws.Range("A57").Value = -22.17430000000001
$ws.Range("B57").Value = 4.693499999999997
$ws.Range("A59").Value = -22.1996
$ws.Range("A69").Value = -21.62219999999999
$ws.Range("B73").Value = 8.223700000000001
$ws.Range("A79").Value = -20.51340000000001
$ws.Range("A83").Value = -21.9761
$ws.Range("B89").Value = 4.949399999999994
$ws.Range("B90").Value = 5.695800000000005
$ws.Range("A93").Value = -21.4019
